$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 191.5
$ws.Range("I9").Value = 191.5
$ws.Range("K9").Value = 191.5
$ws.Range("M9").Value = -22.5

# Row 43
$ws.Range("H43").Value = 1523.1428
$ws.Range("J43").Value = 1541.8462
$ws.Range("L43").Value = 1541.8462
$ws.Range("N43").Value = -1679.8462

# Row 100
$ws.Range("H100").Value = 1051.6
$ws.Range("I100").Value = 1051.6
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1051.6
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -510.5999999999999

# Row 106
$ws.Range("H106").Value = 3901.8572
$ws.Range("I106").Value = 2757.111
$ws.Range("K106").Value = 2757.111
$ws.Range("M106").Value = -2126.111

# Row 132
$ws.Range("H132").Value = 1037.7742
$ws.Range("I132").Value = 1040.3793
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3121.1379
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -591.1379000000002
$ws.Range("N132").Value = -8060

# Row 134
$ws.Range("H134").Value = 44640
$ws.Range("J134").Value = 44640
$ws.Range("L134").Value = 44640
$ws.Range("N134").Value = -54780

# Row 135
$ws.Range("H135").Value = 879.1111
$ws.Range("I135").Value = 879.1111
$ws.Range("K135").Value = 7911.9999
$ws.Range("M135").Value = -5376.9999

# Row 138
$ws.Range("H138").Value = 1830.7667
$ws.Range("I138").Value = 1645.7
$ws.Range("J138").Value = 2200.9
$ws.Range("K138").Value = 4937.1
$ws.Range("L138").Value = 6602.700000000001
$ws.Range("M138").Value = 202.8999999999996
$ws.Range("N138").Value = -16882.7

# Row 141
$ws.Range("H141").Value = 2843.6843
$ws.Range("I141").Value = 1405.6666
$ws.Range("K141").Value = 4216.9998
$ws.Range("M141").Value = 963.0002000000004

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3963.566
$ws.Range("I32").Value = 2144.4146
$ws.Range("K32").Value = 2144.4146
$ws.Range("M32").Value = -1857.4146

# Row 61
$ws.Range("H61").Value = 3728.375
$ws.Range("I61").Value = 2485.3076
$ws.Range("K61").Value = 2485.3076
$ws.Range("M61").Value = -2273.3076

# Row 132
$ws.Range("H132").Value = 1387.8966
$ws.Range("I132").Value = 1102.04
$ws.Range("J132").Value = 3174.5
$ws.Range("K132").Value = 3306.12
$ws.Range("L132").Value = 9523.5
$ws.Range("M132").Value = -776.1199999999999
$ws.Range("N132").Value = -14583.5

# Row 136
$ws.Range("H136").Value = 3728.375
$ws.Range("I136").Value = 2485.3076
$ws.Range("K136").Value = 7455.9228
$ws.Range("M136").Value = -4905.9228

# Row 141
$ws.Range("H141").Value = 29329
$ws.Range("J141").Value = 29329
$ws.Range("L141").Value = 29329
$ws.Range("N141").Value = -39689

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2197.5386
$ws.Range("I20").Value = 2177.9092
$ws.Range("J20").Value = 2305.5
$ws.Range("K20").Value = 2177.9092
$ws.Range("L20").Value = 2305.5
$ws.Range("M20").Value = -1930.9092
$ws.Range("N20").Value = -2799.5

# Row 134
$ws.Range("H134").Value = 12146.192
$ws.Range("I134").Value = 12473.409
$ws.Range("K134").Value = 37420.227
$ws.Range("M134").Value = -34885.227

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 121.875
$ws.Range("I7").Value = 159
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 159
$ws.Range("L7").Value = 60
$ws.Range("M7").Value = -46
$ws.Range("N7").Value = -286

# Row 31
$ws.Range("H31").Value = 3325.0264
$ws.Range("I31").Value = 1616.8667
$ws.Range("K31").Value = 1616.8667
$ws.Range("M31").Value = -1321.8667

# Row 34
$ws.Range("H34").Value = 3325.0264
$ws.Range("I34").Value = 1616.8667
$ws.Range("K34").Value = 1616.8667
$ws.Range("M34").Value = -1414.8667

# Row 58
$ws.Range("H58").Value = 1977746.9
$ws.Range("I58").Value = 3953738.5
$ws.Range("J58").Value = 1755.2727
$ws.Range("K58").Value = 3953738.5
$ws.Range("L58").Value = 1755.2727
$ws.Range("M58").Value = -3953535.5
$ws.Range("N58").Value = -2161.2727

# Row 107
$ws.Range("H107").Value = 335.45
$ws.Range("I107").Value = 335.45
$ws.Range("K107").Value = 335.45
$ws.Range("M107").Value = 1584.55

# Row 134
$ws.Range("H134").Value = 1778.2106
$ws.Range("I134").Value = 1321.5
$ws.Range("K134").Value = 3964.5
$ws.Range("M134").Value = -1429.5

# Row 136
$ws.Range("H136").Value = 1977746.9
$ws.Range("I136").Value = 3953738.5
$ws.Range("J136").Value = 1755.2727
$ws.Range("K136").Value = 11861215.5
$ws.Range("L136").Value = 5265.8181
$ws.Range("M136").Value = -11858665.5
$ws.Range("N136").Value = -10365.8181

# Row 141
$ws.Range("H141").Value = 40096.875
$ws.Range("J141").Value = 40096.875
$ws.Range("L141").Value = 40096.875
$ws.Range("N141").Value = -50456.875

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 1000
$ws.Range("K13").Value = 3000
$ws.Range("M13").Value = -2832

# Row 18
$ws.Range("H18").Value = 1107.25
$ws.Range("I18").Value = 1199.5
$ws.Range("J18").Value = 1015
$ws.Range("K18").Value = 3598.5
$ws.Range("L18").Value = 3045
$ws.Range("M18").Value = -3429.5
$ws.Range("N18").Value = -3383

# Row 131
$ws.Range("H131").Value = 8248.447
$ws.Range("J131").Value = 8972.651
$ws.Range("L131").Value = 26917.953
$ws.Range("N131").Value = -36997.953

# Row 138
$ws.Range("H138").Value = 1805.4
$ws.Range("I138").Value = 1805.4
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 5416.200000000001
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -276.2000000000007

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2985.625
$ws.Range("I102").Value = 3651
$ws.Range("J102").Value = 2468.111
$ws.Range("K102").Value = 3651
$ws.Range("L102").Value = 2468.111
$ws.Range("M102").Value = -2029
$ws.Range("N102").Value = -5712.111

# Row 132
$ws.Range("H132").Value = 1482662.4
$ws.Range("I132").Value = 2138631.8
$ws.Range("J132").Value = 6731.5
$ws.Range("K132").Value = 6415895.399999999
$ws.Range("L132").Value = 20194.5
$ws.Range("M132").Value = -6413365.399999999
$ws.Range("N132").Value = -25254.5

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 1866.7778
$ws.Range("I132").Value = 1998.6666
$ws.Range("J132").Value = 1850.2916
$ws.Range("K132").Value = 5995.9998
$ws.Range("L132").Value = 5550.8748
$ws.Range("M132").Value = -3465.9998
$ws.Range("N132").Value = -10610.8748

# Row 136
$ws.Range("H136").Value = 4038.611
$ws.Range("I136").Value = 2111.2222
$ws.Range("K136").Value = 6333.6666
$ws.Range("M136").Value = -3783.6666

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 7424.893
$ws.Range("I132").Value = 1650
$ws.Range("K132").Value = 4950
$ws.Range("M132").Value = -2420

# Row 136
$ws.Range("H136").Value = 34725540
$ws.Range("I136").Value = 50508310
$ws.Range("K136").Value = 151524930
$ws.Range("M136").Value = -151522380
